$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.354.53'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '1.794.60'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.26'
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.590'
$ws.Range('E6').Value = '  +2.00%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '35.98'
$ws.Range('E8').Value = '  +2.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.290'
$ws.Range('E9').Value = '  -4.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0674'
$ws.Range('E10').Value = '  -3.67%  '
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').Value = '2.052.72'
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.13'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').Value = '1.808.73'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.628'
$ws.Range('E15').Value = '  -2.83%  '
$ws.Range('D16').Value = '34.317.20'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.36'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.66'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.96'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = '0.0₃0767'
$ws.Range('E20').Value = '  -4.13%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.19'
$ws.Range('E22').Value = '  -3.37%  '
$ws.Range('E23').Value = '  -2.55%  '
$ws.Range('E24').Value = '  +3.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.78'
$ws.Range('E25').Value = '  -2.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.87'
$ws.Range('E26').Value = '  +4.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.11'
$ws.Range('E27').Value = '  +1.53%  '
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -1.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.76'
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('E32').Value = '  -3.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0510'
$ws.Range('E33').Value = '  -3.93%  '
$ws.Range('E34').Value = '  -4.69%  '
$ws.Range('D35').Value = '1.355.69'
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.640'
$ws.Range('E36').Value = '  -5.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.05'
$ws.Range('E37').Value = '  -1.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.33'
$ws.Range('E38').Value = '  -9.24%  '
$ws.Range('E39').Value = '  -4.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.41'
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('E41').Value = '  -3.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '80.54'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.926'
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('E44').Value = '  +5.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.05'
$ws.Range('E45').Value = '  -5.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0495'
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('D47').Value = '1.955.25'
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.76'
$ws.Range('E48').Value = '  -4.51%  '
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '101.38'
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -8.87%  '
